# "save data done + era data updated"
# Adds a new "Save" column (H) to the sheet: a header cell formatted like
# the existing headers (bold, bordered, centered - same as G1), plus a
# 0 value in each of the 6 data rows below it. Extends the used range to
# A1:H7 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 - copy G1's formatting (font/border/alignment) so it
# matches the rest of the header row, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data column values (rows 2-7), all 0, default formatting.
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 8).Value = 0
}
